$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("AB2").Value = 1.72

# Row 4
$ws.Range("G4").Value = 1.8
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 4
$ws.Range("N4").Value = 9
$ws.Range("Y4").Value = 1.5
$ws.Range("Z4").Value = 2.5
$ws.Range("AD4").Value = 8
$ws.Range("AQ4").Value = 41

# Row 5
$ws.Range("J5").Value = 2.25
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("AF5").Value = 12
$ws.Range("AJ5").Value = 7
$ws.Range("AQ5").Value = 51
